# Add 2022-Q4 data:
#  - "总计" summary sheet gets a new first-data-row for 2022-Q4 and the
#    previously existing rows shift down (Q3 -> Q3, Q2 -> Q2 positions, i.e.
#    a new row is appended at the bottom for the quarter that is pushed out).
#  - The old "2022-Q3" sheet (holding fund detail numbers) is repurposed to
#    hold the brand-new 2022-Q4 numbers and renamed "2022-Q4".
#  - A fresh sheet named "2022-Q3" is created (an exact duplicate of what the
#    old "2022-Q3" sheet used to contain) so the historical data is preserved.
#  - The "2022-Q2" sheet is left untouched, just shifted to the end.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) overview sheet.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Duplicate row 3 ("2022-Q2", 1, 0.03) down into row 4 so the new row keeps
# the same formatting (bold/bordered index cell in column A) as the others.
$wsTotal.Range("A3:D3").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A4").Value = 2

# Shift the quarter labels down: row2 becomes 2022-Q4 (new), row3 becomes
# 2022-Q3 (previously row2's label), row4 keeps 2022-Q2 (already copied).
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("B2").Value = "2022-Q4"

# ---------------------------------------------------------------------
# 2) Duplicate the existing "2022-Q3" fund-detail sheet so we end up with
#    one copy that keeps the old numbers (renamed "2022-Q3") and one copy
#    that will be overwritten with the new 2022-Q4 numbers (renamed
#    "2022-Q4"), while preserving sheet order 总计, 2022-Q4, 2022-Q3, 2022-Q2.
# ---------------------------------------------------------------------
$wsOldQ3 = $wb.Worksheets.Item("2022-Q3")
$wsOldQ3.Copy($null, $wsOldQ3)
$wsNewQ3 = $wb.Worksheets.Item("2022-Q3 (2)")

$wsOldQ3.Name = "2022-Q4"
$wsNewQ3.Name = "2022-Q3"

$wsQ4 = $wsOldQ3

# ---------------------------------------------------------------------
# 3) Write the new 2022-Q4 numbers on the renamed sheet. Columns D-G hold
#    numeric-looking text (not real numbers) in the source data, so force a
#    text number format before assigning the value and then clear the
#    formatting again so the cell keeps the default (unstyled) appearance.
# ---------------------------------------------------------------------
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $wsQ4.Range("D2") "1.18"
Set-TextValue $wsQ4.Range("E2") "92.77"
Set-TextValue $wsQ4.Range("F2") "2.42"
Set-TextValue $wsQ4.Range("G2") "0.0286"
$wsQ4.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 4) Restore the "2022-Q2" sheet as the active/selected tab, matching the
#    original workbook's view state.
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Activate()
